# Tried to implement Penality Reward System (unfinished)
# Update the "Weekly Quantity" and "Monthly Trend" sheets with revised
# order/requested-quantity figures, and drop rows that are no longer
# part of the analysis window.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The week of 2023-07-24 (serial 45137.99999999999) is removed entirely;
# every later week shifts up by one row.
$ws1.Rows.Item(21).Delete()

# The two most recent weeks (2024-01-07 / 2024-01-14, which after the
# shift above live at rows 27 and 28) are also removed.
$ws1.Rows.Item(27).Delete()
$ws1.Rows.Item(27).Delete()

# Revised requested quantities for the remaining weeks of 2023-09-xx.
$ws1.Range("B18").Value = 2
$ws1.Range("B19").Value = 14
$ws1.Range("B20").Value = 14

# --- Sheet 2: "Monthly Trend" ------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Revised requested quantity for the month of 2023-08 (serial 45138.99999999999).
$ws2.Range("B8").Value = 30

# The most recent month (2024-02, serial 45382.99999999999) is removed.
$ws2.Rows.Item(11).Delete()
